$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.469.38"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "2.639.08"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "2.637.48"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "3.113.09"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "63.389.15"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "2.597.94"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("E25").Value = "  +5.45%  "
$ws.Range("E26").Value = "  +8.23%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "551.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.80%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +12.20%  "
$ws.Range("D34").Value = "0.0₃0809"
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "175.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.52%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").Value = "  +4.01%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "170.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.628"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0552"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0958"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("E51").Value = "  -0.93%  "